$d = $word.ActiveDocument

$pairs = @(
    @("26×33=", "61×17="),
    @("11×77=", "93×50="),
    @("26×25=", "68×71="),
    @("76×36=", "98×55="),
    @("55×68=", "64×25="),
    @("65×32=", "35×25="),
    @("11×21=", "17×18="),
    @("55×30=", "98×53="),
    @("23×73=", "80×35="),
    @("11×49=", "88×76="),
    @("77×75=", "16×92="),
    @("71×41=", "30×49="),
    @("25×77=", "37×60="),
    @("39×27=", "92×80="),
    @("32×52=", "80×83="),
    @("56×92=", "75×11="),
    @("13×32=", "81×94="),
    @("32×80=", "85×42="),
    @("31×51=", "41×69="),
    @("47×69=", "16×94="),
    @("53×98=", "17×64="),
    @("26×86=", "87×20="),
    @("91×86=", "90×21="),
    @("82×73=", "93×35="),
    @("82×82=", "20×57=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
